$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column E (reviews_count) - this shifts F:K left to E:J
$ws.Range("E1").EntireColumn.Delete()
